# Insert a new data row after row 50 (becomes row 51), pushing the existing
# rows 51-105 down to 52-106, then populate the new row with its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(51).Insert()

$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 44874
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 100112021
$ws.Range("G51").Value = "Ají"
$ws.Range("H51").Value = "Americana (o)"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 30
$ws.Range("K51").Value = 35000
$ws.Range("L51").Value = 35000
$ws.Range("M51").Value = 35000
$ws.Range("N51").Value = "`$/caja 25 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 1400
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
